$d = $word.ActiveDocument

function Replace-UniqueText($old, $new) {
    $r = $d.Content
    $r.Start = 0
    $found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: could not find '$old'"
    }
}

# Simple unique numeric / text replacements
Replace-UniqueText "4,103" "2,955"
Replace-UniqueText "786" "729"
Replace-UniqueText "19.2" "24.7"
Replace-UniqueText "562" "560"

# "1,771" must be replaced before the standalone "771" so that the
# substring match inside "1,771" is removed first, leaving "771" unique.
Replace-UniqueText "1,771" "1,746"
Replace-UniqueText "63" "62"
Replace-UniqueText "771" "742"

Replace-UniqueText "939" "938"
Replace-UniqueText "676" "675"
Replace-UniqueText "7,205" "7,148"

Replace-UniqueText "Days, Hauls, Offloads, Trips, Deployments, Marine Mammal Interactions, Samples" "Offloads, Hauls, Trips, Days, Samples, Marine Mammal Interactions, Deployments"
